$d = $word.ActiveDocument

$old = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od Souhvězdí Lva 2022: 14. " + [char]0x2013 + " 23. dubna, 14. " + [char]0x2013 + " 23. května"
$new = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od 14. " + [char]0x2013 + " 23. dubna, 14. " + [char]0x2013 + " 23. května. Při pozorování použijte hvězdy oblohy, které zobrazují souhvězdí Souhvězdí Lva.14. " + [char]0x2013 + " 23. dubna, 14. " + [char]0x2013 + " 23. května"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
